# Append a new scrape run (2026-02-10 07:01:11) to the top of the
# "ランサーズ" sheet, pushing the previous run's surviving rows down,
# dropping the listing that fell off ("Java/講師"), and inserting one
# brand-new listing ("Unity/XRエンジニア募集") in the middle of the
# carried-over rows, matching the scraper's natural write order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# --- column widths (stored OOXML width = ColumnWidth + 5/6, so back it off) ---
$ws.Columns.Item(2).ColumnWidth = 51 - 5/6
$ws.Columns.Item(8).ColumnWidth = 14 - 5/6

# --- row data for the new state (rows 2..10) ---
$timestamp = "2026-02-10 07:01:11"

$data = @(
    @{ B = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集";
       D = "300,000 円 ~ 500,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5423720";
       G = 385;
       H = "🔥AI,Ai ◆効率化" },
    @{ B = "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集";
       D = "200,000 円 ~ 300,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5434128";
       G = 368;
       H = "🔥AI,Ai ◆開発" },
    @{ B = "企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研修講師・メンター)";
       D = "200,000 円 ~ 300,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5434363";
       G = 348;
       H = "🔥AI,Ai ◆コンサル" },
    @{ B = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)";
       D = "300,000 円 ~ 500,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5427956";
       G = 310;
       H = "🔥AI,Ai" },
    @{ B = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪";
       D = "20,000 円 ~ 50,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5217096";
       G = 243;
       H = "🔥API ◆ツール" },
    @{ B = "【急募】新聞記事PDFをCSV・Excel化するPythonプログラム作成依頼";
       D = "50,000 円 ~ 100,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5489128";
       G = 198;
       H = "🔥Python" },
    @{ B = "【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発";
       D = "200,000 円 ~ 300,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5454210";
       G = 108;
       H = "◆開発 ◇アプリ" },
    @{ B = "スプレッドシート(Apps Script)で作業時間をボタン1つで計測・集計できる仕組みの開発";
       D = "50,000 円 ~ 100,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5488743";
       G = 68;
       H = "◆開発" },
    @{ B = "【農機具管理】顧客指定で保有機情報を見れるシステム構築依頼";
       D = "5,000 円 ~ 10,000 円 / 固定";
       F = "https://www.lancers.jp/work/detail/5489112";
       G = 45;
       H = "◇管理" }
)

# Drop every existing hyperlink (the whole sheet's worth); we rebuild them
# below in row order so the relationship ids come out sequential again.
$ws.Range("F2").Hyperlinks.Delete()

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $timestamp
    $ws.Cells.Item($row, 2).Value = $item.B
    $ws.Cells.Item($row, 3).Value = "システム開発"
    $ws.Cells.Item($row, 4).Value = $item.D
    $ws.Cells.Item($row, 5).Value = "期限情報なし"
    $ws.Cells.Item($row, 6).Value = $item.F
    $ws.Cells.Item($row, 7).Value = $item.G
    $ws.Cells.Item($row, 8).Value = $item.H

    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $item.F)

    $row = $row + 1
}
